$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "1a." answer together with
# the hidden _GoBack bookmark at its end ("For this problem we need to
# figure out ... ring finger and so on.").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*For this problem we need to figure out*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the '1a.' finger-counting paragraph"
}

# Split off a new blank paragraph right after it.
$target.Range.InsertParagraphAfter()

# Re-fetch: the paragraph collection shifted, so find the blank paragraph
# we just created (it sits immediately after $target).
$blankIndex = $target.Index + 1
$blank = $d.Paragraphs.Item($blankIndex)

# Insert a second new paragraph after the blank one; this is where the
# new "1b." answer will live.
$blank.Range.InsertParagraphAfter()

$answerIndex = $blankIndex + 1
$answer = $d.Paragraphs.Item($answerIndex)

# Type the new answer text (two sentences, matching the source runs).
$answer.Range.InsertAfter("1b.  As far as insight into this problem, I would say that ")

$answer2 = $d.Paragraphs.Item($answerIndex)
$insertPos = $answer2.Range.End - 1
$tail = $d.Range($insertPos, $insertPos)
$tail.InsertAfter("you should definitely be able to come up with a pattern or an equation to help in solving this, so that you won’t need to just use brute force and count all the way to 1000 on your hand.")

# The _GoBack bookmark originally sat at the end of the "1a." paragraph;
# move it down onto the end of the freshly-typed "1b." paragraph (right
# after the new text, before its paragraph mark) to match the edited
# document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$answer3 = $d.Paragraphs.Item($answerIndex)
$bookmarkPos = $answer3.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
